# Update: Added sales root words
#
# Inserts three new root-word rows into the corpus table on Sheet1:
#   - "money"   under entity "Currency"   (right after the existing "disburse"/Currency row)
#   - "loan"    under entity "Currency"   (right after "money")
#   - "disburse" under entity "Sales"     (right after the existing "supply"/Sales row,
#                                           before the "type"/Product row)
#
# The table has columns: A=id (0-based sequence), B=name (root word), C=root (header only), D=entity
# After the insert the table grows from 67 data rows (A1:D67) to 70 data rows (A1:D70).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the two new "Currency" rows right after the current row 25 (id=23, "disburse") ---
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(26).Insert()

$ws.Cells.Item(26, 2).Value2 = "money"
$ws.Cells.Item(26, 4).Value2 = "Currency"

$ws.Cells.Item(27, 2).Value2 = "loan"
$ws.Cells.Item(27, 4).Value2 = "Currency"

# --- 2. Insert the new "Sales" row (re-using root word "disburse") after the shifted "supply"/Sales row ---
# Before this insert, row 52 = "supply"/Sales and row 53 = "type"/Product (shifted down by the 2 inserts above).
$ws.Rows.Item(53).Insert()

$ws.Cells.Item(53, 2).Value2 = "disburse"
$ws.Cells.Item(53, 4).Value2 = "Sales"

# --- 3. Renumber column A (the 0-based "id" sequence) for every data row, now rows 2..70 ---
$lastRow = 70
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}

# --- 4. Refresh the worksheet AutoFilter so it covers the new extent A1:D70 ---
$full = $ws.Range("A1:D70")
$null = $full.AutoFilter()
$null = $full.AutoFilter()

# --- 5. Update the hidden defined name used by the AutoFilter (_xlnm._FilterDatabase) ---
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -like "*_FilterDatabase*") {
        $nm.RefersTo = "=Sheet1!`$A`$1:`$D`$70"
    }
}
